# Auto update stock data
# Update the report date from 2025/12/28 to 2025/12/29 in column A
# for every stock's first data row (rows 2,8,14,...,74).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    # Keep the cell stored as literal text "2025/12/29" (not an Excel date value),
    # matching how the date was originally stored in the workbook.
    $cell.NumberFormat = "@"
    $cell.Value = "2025/12/29"
}
